$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns holding numeric-looking price strings must stay text so exact
# formatting (trailing zeros, etc.) from the source feed survives the write.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values scraped by the GitHub Actions job.
$ws.Range("D2").Value = "67.840.67"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "2.533.70"
$ws.Range("E3").Value = "  -1.86%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "592.23"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").Value = "175.75"
$ws.Range("E6").Value = "  +3.96%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").Value = "2.531.44"
$ws.Range("E9").Value = "  -1.90%  "
$ws.Range("E10").Value = "  +1.27%  "
$ws.Range("E11").Value = "  +2.20%  "
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "26.79"
$ws.Range("E14").Value = "  -0.32%  "
$ws.Range("D15").Value = "2.993.41"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("E16").Value = "  +0.19%  "
$ws.Range("D17").Value = "67.681.28"
$ws.Range("E17").Value = "  +1.65%  "
$ws.Range("D18").Value = "2.525.14"
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "8.16"
$ws.Range("E19").Value = "  +4.75%  "
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").Value = "359.46"
$ws.Range("E21").Value = "  +2.23%  "
$ws.Range("E22").Value = "  -0.93%  "
$ws.Range("D23").Value = "4.66"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  +4.83%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "10.24"
$ws.Range("E26").Value = "  +3.23%  "
$ws.Range("D27").Value = "70.10"
$ws.Range("E27").Value = "  +1.12%  "
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "0.0₃0994"
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "552.59"
$ws.Range("E31").Value = "  +3.72%  "
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").Value = "1.36"
$ws.Range("E33").Value = "  +1.58%  "
$ws.Range("E34").Value = "  +0.79%  "
$ws.Range("E35").Value = "  -1.38%  "
$ws.Range("D36").Value = "0.998"
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").Value = "157.94"
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").Value = "18.80"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").Value = "18.60"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("B41").Value = "PolygonEcosystemToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D41").Value = "0.356"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "1.82"
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "2.57"
$ws.Range("E44").Value = "  +5.35%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "149.57"
$ws.Range("E46").Value = "  -0.15%  "
$ws.Range("E47").Value = "  -1.15%  "
$ws.Range("D48").Value = "0.0₆0279"
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -1.44%  "
$ws.Range("D51").Value = "0.0759"
$ws.Range("E51").Value = "  -0.57%  "
